# Generate Report for Handoff
#
# Status moved from "In Translation" to "Ready for handoff" on every sheet,
# and the handoff timestamps were refreshed (HO Xliff generate date / latest
# handoff datetime). The "Status" columns also grew wider to fit the new,
# longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps ---------------------------------------------------
# Overview!G2 and de-de!H2 share the "Latest HO Xliff Generate Date" / de-de
# handoff instant.
$wsOverview.Range("G2").Value = "2016-08-22 13:01:46"
$wsDeDe.Range("H2").Value = "2016-08-22 13:01:46"
# zh-cn!H2 is the zh-cn handoff instant.
$wsZhCn.Range("H2").Value = "2016-08-22 13:01:41"

# --- Widen the "Status" columns to fit "Ready for handoff" -----------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
